$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.745.45'
$ws.Range('E2').Value = '  -1.38%  '
$ws.Range('D3').Value = '1.806.57'
$ws.Range('E3').Value = '  -0.92%  '
$ws.Range('E4').Value = '  +0.42%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.55'
$ws.Range('E5').Value = '  -0.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5908'
$ws.Range('E6').Value = '  -1.92%  '
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2777'
$ws.Range('E8').Value = '  -0.46%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06833'
$ws.Range('E9').Value = '  -3.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.31'
$ws.Range('E10').Value = '  -0.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07497'
$ws.Range('E11').Value = '  -1.48%  '
$ws.Range('D12').Value = '1.807.27'
$ws.Range('E12').Value = '  -0.96%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.763'
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6230'
$ws.Range('E14').Value = '  -0.87%  '
$ws.Range('D15').Value = '2.051.90'
$ws.Range('E15').Value = '  -0.88%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009337'
$ws.Range('E16').Value = '  -6.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '75.66'
$ws.Range('E17').Value = '  -3.83%  '
$ws.Range('D18').Value = '28.718.77'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.474'
$ws.Range('E19').Value = '  -6.32%  '
$ws.Range('E20').Value = '  +0.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '211.40'
$ws.Range('E21').Value = '  -6.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.48'
$ws.Range('E22').Value = '  -1.92%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.826'
$ws.Range('E23').Value = '  -2.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.004'
$ws.Range('E24').Value = '  +0.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.98'
$ws.Range('E25').Value = '  -0.64%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.877'
$ws.Range('E26').Value = '  -1.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1268'
$ws.Range('E27').Value = '  -2.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.43'
$ws.Range('E28').Value = '  -0.66%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.439'
$ws.Range('E29').Value = '  -3.37%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.06163'
$ws.Range('E30').Value = '  -0.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.427'
$ws.Range('E31').Value = '  -1.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.786'
$ws.Range('E32').Value = '  -1.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.750'
$ws.Range('E33').Value = '  -0.97%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.735'
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.062'
$ws.Range('E35').Value = '  -5.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6425'
$ws.Range('E36').Value = '  +0.72%  '
$ws.Range('E37').Value = '  -1.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.719'
$ws.Range('E38').Value = '  -0.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.554'
$ws.Range('E39').Value = '  +1.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01699'
$ws.Range('E40').Value = '  -1.83%  '
$ws.Range('D41').Value = '1.147.98'
$ws.Range('E41').Value = '  -5.21%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8828'
$ws.Range('E42').Value = '  -2.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.005'
$ws.Range('E43').Value = '  +0.67%  '
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').Value = '1.959.77'
$ws.Range('E45').Value = '  -1.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '60.52'
$ws.Range('E46').Value = '  -2.99%  '
$ws.Range('E47').Value = '  -2.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.600'
$ws.Range('E48').Value = '  +0.28%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.390'
$ws.Range('E49').Value = '  -1.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05465'
$ws.Range('E50').Value = '  -0.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4477'
$ws.Range('E51').Value = '  -1.62%  '
